$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: update "Actual Result" / "Pass/Fail" for TC_01, and reselect H2:J2 ---
$ws.Range("H2").Value = "Successfully Navigated to Home page"
$ws.Range("I2").Value = "Pass"
$ws.Range("H2:J2").Select()

# --- Row 3: fill in TC_02 details (a new "invalid sign up" test case) ---
$ws.Rows.Item(3).RowHeight = 33.75
$ws.Range("B3").Value = "Invalid Sign up using already registered email"
$ws.Range("C3").Value = "Dummy"
$ws.Range("D3").Value = "Test"
$ws.Range("E3").Value = "CorrectEmail@gmail.com"
$ws.Range("F3").Value = "thatshouldn'tmakeit"

# Add the mailto hyperlink on E3, matching the existing one on E2.
$ws.Hyperlinks.Add($ws.Range("E3"), "mailto:CorrectEmail@gmail.com")
